# Initial Setup of Symptom Elicitation Phase Fulfillment.
#
# Adds a "symptom" training-phrase block to the existing
# `initial_symptom_set` intent row (H17/J17), then appends two new
# intent rows - `fallback_symptom_set` (row 18) and `elicitation`
# (row 19) - on the "intents-en" sheet, replacing the old blank
# placeholder rows 19-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("intents-en")
$ws.Activate()

# --- Remove the old empty placeholder rows (19, 20, 21) ---
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(19).Delete()
$ws.Rows.Item(19).Delete()

# --- Row 18 (new): fallback_symptom_set ---
$ws.Cells.Item(18, 1).Value = "fallback_symptom_set"
$ws.Cells.Item(18, 2).Value = "en"
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = "INITIAL"
$ws.Cells.Item(18, 9).Value = "What symptom?|Could you share one symptom?"
$ws.Rows.Item(18).RowHeight = 22.5

# --- Row 19 (new): elicitation ---
$ws.Cells.Item(19, 1).Value = "elicitation"
$ws.Cells.Item(19, 2).Value = "en"
$ws.Cells.Item(19, 3).Value = 1
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = "ELICITATION"
$ws.Cells.Item(19, 7).Value = "ELICITATION"

# --- Row 17 (initial_symptom_set): fill in training phrases + parameters ---
$ws.Cells.Item(17, 8).Value = "{@symptom:cough}|Have a {@symptom:cough}|Got a {@symptom:cough}|Suffering from {@symptom:cough}|Feeling a {@symptom:cough}|This persistent {@symptom:cough}|This {@symptom:cough}|Dealing with {@symptom:cough}|Down with {@symptom:cough}|Struggling with a {@symptom:cough}|Contracted a {@symptom:cough}|Constant {@symptom:cough}"

$ws.Cells.Item(17, 10).Value = "{@symptom:`$symptom:1:0}"
$ws.Cells.Item(17, 10).HorizontalAlignment = -4131   # xlHAlignLeft
$ws.Cells.Item(17, 10).WrapText = $true

$ws.Rows.Item(17).RowHeight = 135

# --- View state: scroll/zoom/selection to match the edited workbook ---
$win = $ws.Application.ActiveWindow
$win.Zoom = 160
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("G16").Select()
